$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "5" (sheet1.xml): append QAOA4 rows 21-23 to the right-hand table
# (columns M:W). Shared-string order must match: QAOA4(7), QAOA4(9),
# QAOA4(11) were interned in that order, so write row 23 first, then 22,
# then 21, to reproduce the same shared-string insertion order.
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("5")

# Row 23 - QAOA4(7)
$ws5.Cells.Item(23, 13).Value = "QAOA4(7)"
$ws5.Cells.Item(23, 14).Value = 74
$ws5.Cells.Item(23, 15).Value = 63
$ws5.Cells.Item(23, 16).Formula = "=N23-O23"
$ws5.Cells.Item(23, 17).Value = 63
$ws5.Cells.Item(23, 18).Value = 66
$ws5.Cells.Item(23, 19).Formula = "=(O23-Q23)/N23"
$ws5.Cells.Item(23, 20).Formula = "=(P23-R23+Q23)/N23"
$ws5.Cells.Item(23, 21).Formula = "= 1 -R23/N23"
$ws5.Cells.Item(23, 22).Value = 0
$ws5.Cells.Item(23, 23).Formula = "=V23/N23"

# Row 22 - QAOA4(9)
$ws5.Cells.Item(22, 13).Value = "QAOA4(9)"
$ws5.Cells.Item(22, 14).Value = 74
$ws5.Cells.Item(22, 15).Value = 63
$ws5.Cells.Item(22, 16).Formula = "=N22-O22"
$ws5.Cells.Item(22, 17).Value = 40
$ws5.Cells.Item(22, 18).Value = 45
$ws5.Cells.Item(22, 19).Formula = "=(O22-Q22)/N22"
$ws5.Cells.Item(22, 20).Formula = "=(P22-R22+Q22)/N22"
$ws5.Cells.Item(22, 21).Formula = "= 1 -R22/N22"
$ws5.Cells.Item(22, 22).Value = 0
$ws5.Cells.Item(22, 23).Formula = "=V22/N22"

# Row 21 - QAOA4(11)
$ws5.Cells.Item(21, 13).Value = "QAOA4(11)"
$ws5.Cells.Item(21, 14).Value = 74
$ws5.Cells.Item(21, 15).Value = 63
$ws5.Cells.Item(21, 16).Formula = "=N21-O21"
$ws5.Cells.Item(21, 17).Value = 36
$ws5.Cells.Item(21, 18).Value = 41
$ws5.Cells.Item(21, 19).Formula = "=(O21-Q21)/N21"
$ws5.Cells.Item(21, 20).Formula = "=(P21-R21+Q21)/N21"
$ws5.Cells.Item(21, 21).Formula = "= 1 -R21/N21"
$ws5.Cells.Item(21, 23).Formula = "=V21/N21"

$ws5.Columns.Item(13).ColumnWidth = 9.5

# ---------------------------------------------------------------------------
# Sheet "Qubit10" (sheet2.xml): append QAOA8(15) rows 18-20, extending the
# existing shared formulas in P, S, T, U, W downward (these already grow via
# the shared group so a plain copy of the formula cell reproduces it).
# ---------------------------------------------------------------------------
$wsQ = $wb.Worksheets.Item("Qubit10")

# Row 18
$wsQ.Cells.Item(18, 13).Value = "QAOA8(15)"
$wsQ.Cells.Item(18, 14).Value = 120
$wsQ.Cells.Item(18, 15).Value = 109
$wsQ.Cells.Item(18, 16).Formula = "=N18-O18"
$wsQ.Cells.Item(18, 17).Value = 82
$wsQ.Cells.Item(18, 18).Value = 87
$wsQ.Cells.Item(18, 19).Formula = "=(O18-Q18)/N18"
$wsQ.Cells.Item(18, 20).Formula = "=(P18-R18+Q18)/N18"
$wsQ.Cells.Item(18, 21).Formula = "= 1 -R18/N18"
$wsQ.Cells.Item(18, 23).Formula = "=V18/N18"

# Row 19
$wsQ.Cells.Item(19, 13).Value = "QAOA8(15)"
$wsQ.Cells.Item(19, 14).Value = 120
$wsQ.Cells.Item(19, 15).Value = 109
$wsQ.Cells.Item(19, 16).Formula = "=N19-O19"
$wsQ.Cells.Item(19, 17).Value = 92
$wsQ.Cells.Item(19, 18).Value = 98
$wsQ.Cells.Item(19, 19).Formula = "=(O19-Q19)/N19"
$wsQ.Cells.Item(19, 20).Formula = "=(P19-R19+Q19)/N19"
$wsQ.Cells.Item(19, 21).Formula = "= 1 -R19/N19"
$wsQ.Cells.Item(19, 23).Formula = "=V19/N19"

# Row 20
$wsQ.Cells.Item(20, 13).Value = "QAOA8(15)"
$wsQ.Cells.Item(20, 14).Value = 120
$wsQ.Cells.Item(20, 15).Value = 109
$wsQ.Cells.Item(20, 16).Formula = "=N20-O20"
$wsQ.Cells.Item(20, 17).Value = 109
$wsQ.Cells.Item(20, 18).Value = 111
$wsQ.Cells.Item(20, 19).Formula = "=(O20-Q20)/N20"
$wsQ.Cells.Item(20, 20).Formula = "=(P20-R20+Q20)/N20"
$wsQ.Cells.Item(20, 21).Formula = "= 1 -R20/N20"
$wsQ.Cells.Item(20, 22).Value = 0
$wsQ.Cells.Item(20, 23).Formula = "=V20/N20"

$wsQ.Columns.Item(13).ColumnWidth = 11.5

# ---------------------------------------------------------------------------
# Selections: restore the recorded cursor position on each sheet, and make
# sure sheet "5" ends up the active tab (matches the saved file).
# ---------------------------------------------------------------------------
$wsQ.Activate()
$wsQ.Range("U30").Select()

$ws5.Activate()
$ws5.Range("R21").Select()
